$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 84: correct the timestamp in column A ---
$ws.Cells.Item(84, 1).Value = 45463.2916666667

# --- Append new row 85 with the new data point ---
$ws.Cells.Item(84, 1).Copy()
$ws.Cells.Item(85, 1).PasteSpecial(-4122)
$ws.Cells.Item(85, 1).Value = 45464.3175462963

$ws.Cells.Item(85, 2).Value = 1500
$ws.Cells.Item(85, 3).Value = 2.95000004768372
$ws.Cells.Item(85, 4).Value = 2.95000004768372
$ws.Cells.Item(85, 5).Value = 2.95000004768372
$ws.Cells.Item(85, 6).Value = 2.95000004768372
$ws.Cells.Item(85, 7).NumberFormat = "@"
$ws.Cells.Item(85, 7).Value = "2.95000004768372"
$ws.Cells.Item(85, 7).ClearFormats()
$ws.Cells.Item(85, 8).Value = "ESPE.MI"
